$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Cadastrado" (B) and "Sem Cadastro" (C) values for each year row
$ws.Range("B2").Value = 101964.37
$ws.Range("C2").Value = 13544.65

$ws.Range("C3").Value = 46444.4

$ws.Range("B4").Value = 1749646.5
$ws.Range("C4").Value = 22553.9

$ws.Range("B5").Value = 2822961.3
$ws.Range("C5").Value = 20181.14

$ws.Range("B6").Value = 4416713.12
$ws.Range("C6").Value = 26365

$ws.Range("B7").Value = 2268166.05
$ws.Range("C7").Value = 12953.54
